$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 710.375
$ws.Range("I33").Value = 581
$ws.Range("K33").Value = 581
$ws.Range("M33").Value = -352

$ws.Range("H98").Value = 2170.75
$ws.Range("I98").Value = 1947.0857
$ws.Range("K98").Value = 1947.0857
$ws.Range("M98").Value = -449.0857000000001

$ws.Range("H122").Value = 2170.75
$ws.Range("I122").Value = 1947.0857
$ws.Range("K122").Value = 5841.257100000001
$ws.Range("M122").Value = -3391.257100000001

$ws.Range("H125").Value = 2016.3334
$ws.Range("J125").Value = 2088.3333
$ws.Range("L125").Value = 18794.9997
$ws.Range("N125").Value = -23714.9997

$ws.Range("H127").Value = 992.4666999999999
$ws.Range("I127").Value = 1297.6666
$ws.Range("J127").Value = 958.55554
$ws.Range("K127").Value = 3892.9998
$ws.Range("L127").Value = 2875.66662
$ws.Range("M127").Value = 1067.0002
$ws.Range("N127").Value = -12795.66662

$ws.Range("H129").Value = 3572371
$ws.Range("I129").Value = 50001420
$ws.Range("J129").Value = 905.6923
$ws.Range("K129").Value = 150004260
$ws.Range("L129").Value = 2717.0769
$ws.Range("M129").Value = -149999260
$ws.Range("N129").Value = -12717.0769

$ws.Range("H138").Value = 5228.778
$ws.Range("I138").Value = 2771.7
$ws.Range("J138").Value = 6034.377
$ws.Range("K138").Value = 8315.099999999999
$ws.Range("L138").Value = 18103.131
$ws.Range("M138").Value = -3175.099999999999
$ws.Range("N138").Value = -28383.131

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2157.6
$ws.Range("I132").Value = 1759.2759
$ws.Range("J132").Value = 4082.8333
$ws.Range("K132").Value = 5277.8277
$ws.Range("L132").Value = 12248.4999
$ws.Range("M132").Value = -2747.8277
$ws.Range("N132").Value = -17308.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 41061.54
$ws.Range("J86").Value = 44366.668
$ws.Range("L86").Value = 44366.668
$ws.Range("N86").Value = -46612.668

$ws.Range("H89").Value = 41061.54
$ws.Range("J89").Value = 44366.668
$ws.Range("L89").Value = 221833.34
$ws.Range("N89").Value = -233065.34

$ws.Range("H99").Value = 3136.125
$ws.Range("I99").Value = 2629.077
$ws.Range("J99").Value = 5333.3335
$ws.Range("K99").Value = 2629.077
$ws.Range("L99").Value = 5333.3335
$ws.Range("M99").Value = -1131.077
$ws.Range("N99").Value = -8329.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3256.2778
$ws.Range("I31").Value = 2217.8684
$ws.Range("J31").Value = 4416.853
$ws.Range("K31").Value = 2217.8684
$ws.Range("L31").Value = 4416.853
$ws.Range("M31").Value = -1922.8684
$ws.Range("N31").Value = -5006.853

$ws.Range("H34").Value = 3256.2778
$ws.Range("I34").Value = 2217.8684
$ws.Range("J34").Value = 4416.853
$ws.Range("K34").Value = 2217.8684
$ws.Range("L34").Value = 4416.853
$ws.Range("M34").Value = -2015.8684
$ws.Range("N34").Value = -4820.853

$ws.Range("H132").Value = 3628.5264
$ws.Range("I132").Value = 3253.1428
$ws.Range("J132").Value = 4679.6
$ws.Range("K132").Value = 9759.428400000001
$ws.Range("L132").Value = 14038.8
$ws.Range("M132").Value = -7229.428400000001
$ws.Range("N132").Value = -19098.8

$ws.Range("H134").Value = 16132284
$ws.Range("I134").Value = 25003040
$ws.Range("J134").Value = 3635.7273
$ws.Range("K134").Value = 75009120
$ws.Range("L134").Value = 10907.1819
$ws.Range("M134").Value = -75006585
$ws.Range("N134").Value = -15977.1819

$ws.Range("H141").Value = 26106.666
$ws.Range("J141").Value = 26106.666
$ws.Range("L141").Value = 26106.666
$ws.Range("N141").Value = -36466.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5531.4614
$ws.Range("I56").Value = 5531.4614
$ws.Range("K56").Value = 5531.4614
$ws.Range("M56").Value = -5001.4614

$ws.Range("H107").Value = 1327.4706
$ws.Range("I107").Value = 376.3
$ws.Range("J107").Value = 2686.2856
$ws.Range("K107").Value = 1128.9
$ws.Range("L107").Value = 8058.8568
$ws.Range("M107").Value = 791.0999999999999
$ws.Range("N107").Value = -11898.8568

$ws.Range("H112").Value = 3526.6667
$ws.Range("I112").Value = 1800
$ws.Range("J112").Value = 3683.6365
$ws.Range("K112").Value = 5400
$ws.Range("L112").Value = 11050.9095
$ws.Range("M112").Value = -4292
$ws.Range("N112").Value = -13266.9095

$ws.Range("H113").Value = 882.3103599999999
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 895.9643
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2687.8929
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -7027.8929

$ws.Range("H119").Value = 3006.8462
$ws.Range("I119").Value = 1417.8
$ws.Range("K119").Value = 4253.4
$ws.Range("M119").Value = 584.6000000000004

$ws.Range("H131").Value = 1377.1746
$ws.Range("I131").Value = 3628.5715
$ws.Range("J131").Value = 1095.75
$ws.Range("K131").Value = 10885.7145
$ws.Range("L131").Value = 3287.25
$ws.Range("M131").Value = -5845.7145
$ws.Range("N131").Value = -13367.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4899.2856
$ws.Range("I70").Value = 4648.5713
$ws.Range("J70").Value = 5150
$ws.Range("K70").Value = 4648.5713
$ws.Range("L70").Value = 5150
$ws.Range("M70").Value = -4378.5713
$ws.Range("N70").Value = -5690

$ws.Range("H73").Value = 4899.2856
$ws.Range("I73").Value = 4648.5713
$ws.Range("J73").Value = 5150
$ws.Range("K73").Value = 4648.5713
$ws.Range("L73").Value = 5150
$ws.Range("M73").Value = -3712.5713
$ws.Range("N73").Value = -7022

$ws.Range("H126").Value = 838034.5
$ws.Range("I126").Value = 4100
$ws.Range("J126").Value = 1255001.8
$ws.Range("K126").Value = 12300
$ws.Range("L126").Value = 3765005.4
$ws.Range("M126").Value = -9830
$ws.Range("N126").Value = -3769945.4

$ws.Range("H132").Value = 4737.1353
$ws.Range("I132").Value = 5926.4443
$ws.Range("J132").Value = 3610.4211
$ws.Range("K132").Value = 17779.3329
$ws.Range("L132").Value = 10831.2633
$ws.Range("M132").Value = -15249.3329
$ws.Range("N132").Value = -15891.2633

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 29966.666
$ws.Range("J141").Value = 29966.666
$ws.Range("L141").Value = 29966.666
$ws.Range("N141").Value = -40326.666

Write-Host "Applied all scheduled Sheets updates"
